$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the generic-name / manufacturer cells that become blank for these rows
# (mirrors rows where the values were removed in the target workbook)
$ws.Range("B3").ClearContents()
$ws.Range("B6").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("B9").ClearContents()
$ws.Range("H5").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("H9").ClearContents()

# Recalculate so the concatenation formulas in column A reflect the cleared cells
$excel.Calculate()

# Update the active selection to H10
$ws.Range("H10").Select()
